$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.556.66'
$ws.Range("E2").Value = '  +0.97%  '
$ws.Range("D3").Value = '2.499.45'
$ws.Range("E3").Value = '  +2.17%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.992'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.77%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.64'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.26%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '94.76'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.77%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.548'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.38%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.994'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.72%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.501'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.28%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '33.07'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.47%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0787'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.05%  '
$ws.Range("E12").Value = '  +1.90%  '
$ws.Range("D13").Value = '2.876.42'
$ws.Range("E13").Value = '  +2.00%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.91'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.52%  '
$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '2.542.22'
$ws.Range("E15").Value = '  +3.81%  '
$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.59'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +8.34%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.763'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.47%  '
$ws.Range("D18").Value = '41.701.60'
$ws.Range("E18").Value = '  +1.40%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.34'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.36%  '
$ws.Range("D20").Value = '0.0₃0926'
$ws.Range("E20").Value = '  +1.88%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '70.71'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.83%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.29'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.62%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.50'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.05%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.72'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.47%  '
$ws.Range("B25").Value = 'ImmutableX'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.91'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.91%  '
$ws.Range("B26").Value = 'Dai'
$ws.Range("C26").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.10%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.81'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.45%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.24'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.33%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.72'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.69%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.57'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.76%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '154.69'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.83%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.45'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.32%  '
$ws.Range("B33").Value = 'Celestia'
$ws.Range("C33").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.44'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +7.76%  '
$ws.Range("B34").Value = 'WEMIXToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.58'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.21%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0762'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.57%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.49'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.29%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.00'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.37%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.85'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.43%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.114'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.11%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.102'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.70%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.16'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.28%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.995'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.88%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '19.83'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.23%  '
$ws.Range("D44").Value = '1.960.62'
$ws.Range("E44").Value = '  -0.34%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0286'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.10%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.00'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.22%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.85'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.67%  '
$ws.Range("D48").Value = '2.733.72'
$ws.Range("E48").Value = '  +2.02%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '96.88'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.07%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '67.94'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.13%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.177'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.31%  '
